$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.112.77"
$ws.Range("E2").Value = "  -0.38%  "

Set-TextValue $ws.Range("D3") "1.827.97"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("E4").Value = "  -0.38%  "

Set-TextValue $ws.Range("D5") "311.60"
$ws.Range("E5").Value = "  -0.91%  "

Set-TextValue $ws.Range("D6") "0.9996"
$ws.Range("E6").Value = "  -0.40%  "

Set-TextValue $ws.Range("D7") "0.5107"
$ws.Range("E7").Value = "  -2.84%  "

Set-TextValue $ws.Range("D8") "0.3960"
$ws.Range("E8").Value = "  +3.73%  "

Set-TextValue $ws.Range("D9") "0.1011"
$ws.Range("E9").Value = "  +26.70%  "

Set-TextValue $ws.Range("D10") "1.111"
$ws.Range("E10").Value = "  +1.08%  "

Set-TextValue $ws.Range("D11") "41.00"
$ws.Range("E11").Value = "  -0.75%  "

Set-TextValue $ws.Range("D12") "6.481"
$ws.Range("E12").Value = "  +2.73%  "

Set-TextValue $ws.Range("D13") "1.000"
$ws.Range("E13").Value = "  -0.43%  "

Set-TextValue $ws.Range("D14") "20.64"
$ws.Range("E14").Value = "  +0.31%  "

Set-TextValue $ws.Range("D15") "7.401"
$ws.Range("E15").Value = "  +1.21%  "

Set-TextValue $ws.Range("D16") "1.822.60"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "95.30"
$ws.Range("E17").Value = "  +3.47%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001139"
$ws.Range("E18").Value = "  +4.15%  "

Set-TextValue $ws.Range("D19") "0.06606"
$ws.Range("E19").Value = "  +0.14%  "

Set-TextValue $ws.Range("D20") "0.9997"
$ws.Range("E20").Value = "  -0.39%  "

Set-TextValue $ws.Range("D21") "17.36"
$ws.Range("E21").Value = "  +0.09%  "

Set-TextValue $ws.Range("D22") "6.054"
$ws.Range("E22").Value = "  +1.46%  "

Set-TextValue $ws.Range("D23") "28.190.99"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("E24").Value = "  +0.57%  "

Set-TextValue $ws.Range("D25") "2.243"
$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D26") "2.477"
$ws.Range("E26").Value = "  +5.20%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "158.61"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "20.80"
$ws.Range("E28").Value = "  +1.74%  "

Set-TextValue $ws.Range("D29") "2.035.41"
$ws.Range("E29").Value = "  +1.28%  "

Set-TextValue $ws.Range("D30") "128.62"
$ws.Range("E30").Value = "  +4.49%  "

Set-TextValue $ws.Range("D31") "0.1096"
$ws.Range("E31").Value = "  +1.25%  "

Set-TextValue $ws.Range("D32") "1.066"
$ws.Range("E32").Value = "  +0.95%  "

Set-TextValue $ws.Range("D33") "5.645"
$ws.Range("E33").Value = "  +1.75%  "

Set-TextValue $ws.Range("D34") "3.638"
$ws.Range("E34").Value = "  -1.36%  "

Set-TextValue $ws.Range("D35") "0.06904"
$ws.Range("E35").Value = "  -4.29%  "

Set-TextValue $ws.Range("D36") "9.178"
$ws.Range("E36").Value = "  +6.68%  "

Set-TextValue $ws.Range("D37") "0.02346"
$ws.Range("E37").Value = "  +1.61%  "

Set-TextValue $ws.Range("D38") "0.2171"
$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("E39").Value = "  -5.91%  "

Set-TextValue $ws.Range("D40") "5.028"
$ws.Range("E40").Value = "  -1.25%  "

Set-TextValue $ws.Range("D41") "0.6270"
$ws.Range("E41").Value = "  +1.32%  "

Set-TextValue $ws.Range("D42") "0.9998"
$ws.Range("E42").Value = "  -0.22%  "

Set-TextValue $ws.Range("D43") "1.159"
$ws.Range("E43").Value = "  -0.72%  "

Set-TextValue $ws.Range("D44") "13.32"
$ws.Range("E44").Value = "  +0.88%  "

Set-TextValue $ws.Range("D45") "0.6002"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D46") "3.709"
$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D47") "1.287"
$ws.Range("E47").Value = "  -6.22%  "

Set-TextValue $ws.Range("D48") "125.86"
$ws.Range("E48").Value = "  -1.05%  "

Set-TextValue $ws.Range("D49") "1.992"
$ws.Range("E49").Value = "  +3.57%  "

$ws.Range("E50").Value = "  -2.52%  "

Set-TextValue $ws.Range("D51") "0.06788"
$ws.Range("E51").Value = "  -0.27%  "
